$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.223.13"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").Value = "3.426.94"
$ws.Range("E3").Value = "  +0.27%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.29%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "413.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.14%  "

$ws.Range("E7").Value = "  -1.17%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.724"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.00%  "

$ws.Range("E10").Value = "  +1.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.63"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.90%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.94%  "

$ws.Range("E13").Value = "  +9.13%  "

$ws.Range("D14").Value = "3.973.86"
$ws.Range("E14").Value = "  +0.24%  "

$ws.Range("E15").Value = "  -0.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.64%  "

$ws.Range("D17").Value = "3.427.06"
$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.49%  "

$ws.Range("E19").Value = "  +0.13%  "

$ws.Range("D20").Value = "62.273.44"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "471.31"
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "90.98"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +21.00%  "

$ws.Range("E26").Value = "  +2.29%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.12%  "

$ws.Range("E28").Value = "  -0.38%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.94%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.72%  "

$ws.Range("E31").Value = "  -3.85%  "

$ws.Range("E32").Value = "  -2.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.112"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.80%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.59"
$ws.Range("D34").Style = "Normal"

$ws.Range("E35").Value = "  +0.08%  "

$ws.Range("E36").Value = "  +10.84%  "

$ws.Range("E37").Value = "  -1.52%  "

$ws.Range("E38").Value = "  +0.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.22%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.324"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.10%  "

$ws.Range("E41").Value = "  +0.58%  "

$ws.Range("E42").Value = "  -1.08%  "

$ws.Range("E43").Value = "  +11.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "145.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.58%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +17.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "16.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.87%  "

$ws.Range("D49").Value = "0.0₃0538"
$ws.Range("E49").Value = "  +36.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.70%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "110.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.66%  "
